$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (old row 2); this shifts subsequent rows up by one
$ws.Rows("2:2").Delete() | Out-Null

# Apply corrected / recomputed values for the forecast data (bugfix)
$ws.Cells.Item(2, 1).Value2 = 39583
$ws.Cells.Item(2, 2).Value2 = 2008
$ws.Cells.Item(2, 3).ClearContents() | Out-Null
$ws.Cells.Item(2, 4).Value2 = 2009
$ws.Cells.Item(2, 5).Value2 = 0.5087393606160395

$ws.Cells.Item(3, 1).Value2 = 39765
$ws.Cells.Item(3, 2).Value2 = 2008
$ws.Cells.Item(3, 3).ClearContents() | Out-Null
$ws.Cells.Item(3, 4).Value2 = 2009
$ws.Cells.Item(3, 5).Value2 = -0.4513776153963867

$ws.Cells.Item(4, 1).Value2 = 39948
$ws.Cells.Item(4, 2).Value2 = 2009
$ws.Cells.Item(4, 3).Value2 = -1.118515468742087
$ws.Cells.Item(4, 4).Value2 = 2010
$ws.Cells.Item(4, 5).Value2 = -0.6296678961043134

$ws.Cells.Item(5, 1).Value2 = 40130
$ws.Cells.Item(5, 2).Value2 = 2009
$ws.Cells.Item(5, 3).Value2 = -1.324983933426882
$ws.Cells.Item(5, 4).Value2 = 2010
$ws.Cells.Item(5, 5).Value2 = -0.8803581938132576

$ws.Cells.Item(6, 1).Value2 = 40310
$ws.Cells.Item(6, 2).Value2 = 2010
$ws.Cells.Item(6, 3).Value2 = -0.1156872058426073
$ws.Cells.Item(6, 4).Value2 = 2011
$ws.Cells.Item(6, 5).Value2 = -0.5120992642018263

$ws.Cells.Item(7, 1).Value2 = 40494
$ws.Cells.Item(7, 2).Value2 = 2010
$ws.Cells.Item(7, 3).Value2 = -0.3900454704678369
$ws.Cells.Item(7, 4).Value2 = 2011
$ws.Cells.Item(7, 5).Value2 = -1.213027585730386

$ws.Cells.Item(8, 1).Value2 = 40676
$ws.Cells.Item(8, 2).Value2 = 2011
$ws.Cells.Item(8, 3).Value2 = -0.4084169314491404
$ws.Cells.Item(8, 4).Value2 = 2012
$ws.Cells.Item(8, 5).Value2 = -0.6403426624573716

$ws.Cells.Item(9, 1).Value2 = 40862
$ws.Cells.Item(9, 2).Value2 = 2011
$ws.Cells.Item(9, 3).Value2 = -0.2995848153489522
$ws.Cells.Item(9, 4).Value2 = 2012
$ws.Cells.Item(9, 5).Value2 = -0.3230872999110068

$ws.Cells.Item(10, 1).Value2 = 41044
$ws.Cells.Item(10, 2).Value2 = 2012
$ws.Cells.Item(10, 3).Value2 = -0.2188016966516937
$ws.Cells.Item(10, 4).Value2 = 2013
$ws.Cells.Item(10, 5).Value2 = -0.1561757764150462

$ws.Cells.Item(11, 1).Value2 = 41228
$ws.Cells.Item(11, 2).Value2 = 2012
$ws.Cells.Item(11, 3).Value2 = -0.2075757021743008
$ws.Cells.Item(11, 4).Value2 = 2013
$ws.Cells.Item(11, 5).Value2 = -0.2793004163246238

$ws.Cells.Item(12, 1).Value2 = 41409
$ws.Cells.Item(12, 2).Value2 = 2013
$ws.Cells.Item(12, 3).Value2 = 0.05915234751026066
$ws.Cells.Item(12, 4).Value2 = 2014
$ws.Cells.Item(12, 5).Value2 = 0.04624521867206965

$ws.Cells.Item(13, 1).Value2 = 41592
$ws.Cells.Item(13, 2).Value2 = 2013
$ws.Cells.Item(13, 3).Value2 = 0.124712275190686
$ws.Cells.Item(13, 4).Value2 = 2014
$ws.Cells.Item(13, 5).Value2 = -0.119752617912039

$ws.Cells.Item(14, 1).Value2 = 41774
$ws.Cells.Item(14, 2).Value2 = 2014
$ws.Cells.Item(14, 3).Value2 = -0.2979029954603529
$ws.Cells.Item(14, 4).Value2 = 2015
$ws.Cells.Item(14, 5).Value2 = -0.1124510725819206

$ws.Cells.Item(15, 1).Value2 = 41957
$ws.Cells.Item(15, 2).Value2 = 2014
$ws.Cells.Item(15, 3).Value2 = -0.255298189276465
$ws.Cells.Item(15, 4).Value2 = 2015
$ws.Cells.Item(15, 5).Value2 = -0.05946205208092747

$ws.Cells.Item(16, 1).Value2 = 42137
$ws.Cells.Item(16, 2).Value2 = 2015
$ws.Cells.Item(16, 3).Value2 = 0.07317408757452348
$ws.Cells.Item(16, 4).Value2 = 2016
$ws.Cells.Item(16, 5).Value2 = -0.002181547367274828

$ws.Cells.Item(17, 1).Value2 = 42321
$ws.Cells.Item(17, 2).Value2 = 2015
$ws.Cells.Item(17, 3).Value2 = 0.07418514192796266
$ws.Cells.Item(17, 4).Value2 = 2016
$ws.Cells.Item(17, 5).Value2 = -0.001680662521774678

$ws.Cells.Item(18, 1).Value2 = 42503
$ws.Cells.Item(18, 2).Value2 = 2016
$ws.Cells.Item(18, 3).Value2 = -0.06188089372189953
$ws.Cells.Item(18, 4).Value2 = 2017
$ws.Cells.Item(18, 5).Value2 = -0.07932008107318644

$ws.Cells.Item(19, 1).Value2 = 42689
$ws.Cells.Item(19, 2).Value2 = 2016
$ws.Cells.Item(19, 3).Value2 = -0.07611406013281474
$ws.Cells.Item(19, 4).Value2 = 2017
$ws.Cells.Item(19, 5).Value2 = -0.1247901924724348

$ws.Cells.Item(20, 1).Value2 = 42867
$ws.Cells.Item(20, 2).Value2 = 2017
$ws.Cells.Item(20, 3).Value2 = -0.2199961235931358
$ws.Cells.Item(20, 4).Value2 = 2018
$ws.Cells.Item(20, 5).Value2 = -0.1022879117640763

$ws.Cells.Item(21, 1).Value2 = 43053
$ws.Cells.Item(21, 2).Value2 = 2017
$ws.Cells.Item(21, 3).Value2 = -0.191300579729714
$ws.Cells.Item(21, 4).Value2 = 2018
$ws.Cells.Item(21, 5).Value2 = -0.05219951976568327

$ws.Cells.Item(22, 1).Value2 = 43145
$ws.Cells.Item(22, 2).Value2 = 2018
$ws.Cells.Item(22, 3).Value2 = 0.0441865668729946
$ws.Cells.Item(22, 4).Value2 = 2019
$ws.Cells.Item(22, 5).Value2 = -0.07120909843567613

$ws.Cells.Item(23, 1).Value2 = 43235
$ws.Cells.Item(23, 2).Value2 = 2018
$ws.Cells.Item(23, 3).Value2 = 0.0882025545300813
$ws.Cells.Item(23, 4).Value2 = 2019
$ws.Cells.Item(23, 5).Value2 = -0.05573300569792217

$ws.Cells.Item(24, 1).Value2 = 43326
$ws.Cells.Item(24, 2).Value2 = 2018
$ws.Cells.Item(24, 3).Value2 = 0.1415113532986956
$ws.Cells.Item(24, 4).Value2 = 2019
$ws.Cells.Item(24, 5).Value2 = 0.02926805735909976

$ws.Cells.Item(25, 1).Value2 = 43418
$ws.Cells.Item(25, 2).Value2 = 2018
$ws.Cells.Item(25, 3).Value2 = 0.0970330232288763
$ws.Cells.Item(25, 4).Value2 = 2019
$ws.Cells.Item(25, 5).Value2 = -0.1345737582127748

$ws.Cells.Item(26, 1).Value2 = 43510
$ws.Cells.Item(26, 2).Value2 = 2019
$ws.Cells.Item(26, 3).Value2 = -0.5756287392657988
$ws.Cells.Item(26, 4).Value2 = 2020
$ws.Cells.Item(26, 5).Value2 = -0.2394607875814136

$ws.Cells.Item(27, 1).Value2 = 43600
$ws.Cells.Item(27, 2).Value2 = 2019
$ws.Cells.Item(27, 3).Value2 = -0.7844010209450802
$ws.Cells.Item(27, 4).Value2 = 2020
$ws.Cells.Item(27, 5).Value2 = -0.3786583343736716

$ws.Cells.Item(28, 1).Value2 = 43691
$ws.Cells.Item(28, 2).Value2 = 2019
$ws.Cells.Item(28, 3).Value2 = -0.6919146680131605
$ws.Cells.Item(28, 4).Value2 = 2020
$ws.Cells.Item(28, 5).Value2 = -0.244860729922769

$ws.Cells.Item(29, 1).Value2 = 43783
$ws.Cells.Item(29, 2).Value2 = 2019
$ws.Cells.Item(29, 3).Value2 = -0.7407518902333265
$ws.Cells.Item(29, 4).Value2 = 2020
$ws.Cells.Item(29, 5).Value2 = -0.4363737508290888

$ws.Cells.Item(30, 1).Value2 = 43875
$ws.Cells.Item(30, 2).Value2 = 2020
$ws.Cells.Item(30, 3).Value2 = -0.7124953797697064
$ws.Cells.Item(30, 4).Value2 = 2021
$ws.Cells.Item(30, 5).Value2 = -0.4617192974095352

$ws.Cells.Item(31, 1).Value2 = 43966
$ws.Cells.Item(31, 2).Value2 = 2020
$ws.Cells.Item(31, 3).Value2 = -0.1808804304865297
$ws.Cells.Item(31, 4).Value2 = 2021
$ws.Cells.Item(31, 5).Value2 = -0.1077309791980285

$ws.Cells.Item(32, 1).Value2 = 44068
$ws.Cells.Item(32, 2).Value2 = 2020
$ws.Cells.Item(32, 3).Value2 = 0.3056679541520335
$ws.Cells.Item(32, 4).Value2 = 2021
$ws.Cells.Item(32, 5).Value2 = 0.3245880452514394

$ws.Cells.Item(33, 1).Value2 = 44159
$ws.Cells.Item(33, 2).Value2 = 2020
$ws.Cells.Item(33, 3).Value2 = 0.3056679541520335
$ws.Cells.Item(33, 4).Value2 = 2021
$ws.Cells.Item(33, 5).Value2 = -0.514812792200714

$ws.Cells.Item(34, 1).Value2 = 44251
$ws.Cells.Item(34, 2).Value2 = 2021
$ws.Cells.Item(34, 3).Value2 = -0.8680533514735522
$ws.Cells.Item(34, 4).Value2 = 2022
$ws.Cells.Item(34, 5).Value2 = -0.5995895195426981

$ws.Cells.Item(35, 1).Value2 = 44341
$ws.Cells.Item(35, 2).Value2 = 2021
$ws.Cells.Item(35, 3).Value2 = -0.8769761459347714
$ws.Cells.Item(35, 4).Value2 = 2022
$ws.Cells.Item(35, 5).Value2 = -0.5354669478056073

$ws.Cells.Item(36, 1).Value2 = 44432
$ws.Cells.Item(36, 2).Value2 = 2021
$ws.Cells.Item(36, 3).Value2 = -1.388491535160907
$ws.Cells.Item(36, 4).Value2 = 2022
$ws.Cells.Item(36, 5).Value2 = -2.541003699199929

$ws.Cells.Item(37, 1).Value2 = 44525
$ws.Cells.Item(37, 2).Value2 = 2021
$ws.Cells.Item(37, 3).Value2 = -1.388491535160907
$ws.Cells.Item(37, 4).Value2 = 2022
$ws.Cells.Item(37, 5).Value2 = -2.321721165370549

$ws.Cells.Item(38, 1).Value2 = 44617
$ws.Cells.Item(38, 2).Value2 = 2022
$ws.Cells.Item(38, 3).Value2 = -1.867377038014506
$ws.Cells.Item(38, 4).Value2 = 2023
$ws.Cells.Item(38, 5).Value2 = -0.8537083312609495

$ws.Cells.Item(39, 1).Value2 = 44706
$ws.Cells.Item(39, 2).Value2 = 2022
$ws.Cells.Item(39, 3).Value2 = -1.719168896439693
$ws.Cells.Item(39, 4).Value2 = 2023
$ws.Cells.Item(39, 5).Value2 = -0.5989817782328322

$ws.Cells.Item(40, 1).Value2 = 44798
$ws.Cells.Item(40, 2).Value2 = 2022
$ws.Cells.Item(40, 3).Value2 = -1.678482969789596
$ws.Cells.Item(40, 4).Value2 = 2023
$ws.Cells.Item(40, 5).Value2 = -0.6229862770763095

$ws.Cells.Item(41, 1).Value2 = 44890
$ws.Cells.Item(41, 2).Value2 = 2022
$ws.Cells.Item(41, 3).Value2 = -1.678482969789596
$ws.Cells.Item(41, 4).Value2 = 2023
$ws.Cells.Item(41, 5).Value2 = -1.107351089172237

$ws.Cells.Item(42, 1).Value2 = 44981
$ws.Cells.Item(42, 2).Value2 = 2023
$ws.Cells.Item(42, 3).Value2 = -0.9016470784766528
$ws.Cells.Item(42, 4).Value2 = 2024
$ws.Cells.Item(42, 5).Value2 = -1.556833564431637

$ws.Cells.Item(43, 1).Value2 = 45071
$ws.Cells.Item(43, 2).Value2 = 2023
$ws.Cells.Item(43, 3).Value2 = -0.6902657121583777
$ws.Cells.Item(43, 4).Value2 = 2024
$ws.Cells.Item(43, 5).Value2 = -1.134712300966823

$ws.Cells.Item(44, 1).Value2 = 45163
$ws.Cells.Item(44, 2).Value2 = 2023
$ws.Cells.Item(44, 3).Value2 = -0.5999457276250508
$ws.Cells.Item(44, 4).Value2 = 2024
$ws.Cells.Item(44, 5).Value2 = -0.7253995615808195

$ws.Cells.Item(45, 1).Value2 = 45254
$ws.Cells.Item(45, 2).Value2 = 2023
$ws.Cells.Item(45, 3).Value2 = -0.5999457276250508
$ws.Cells.Item(45, 4).Value2 = 2024
$ws.Cells.Item(45, 5).Value2 = -0.4628630633218611

$ws.Cells.Item(46, 1).Value2 = 45345
$ws.Cells.Item(46, 2).Value2 = 2024
$ws.Cells.Item(46, 3).Value2 = -0.101315145211045
$ws.Cells.Item(46, 4).Value2 = 2025
$ws.Cells.Item(46, 5).Value2 = -0.1859728711464226

$ws.Cells.Item(47, 1).Value2 = 45436
$ws.Cells.Item(47, 2).Value2 = 2024
$ws.Cells.Item(47, 3).Value2 = -0.07874066250703748
$ws.Cells.Item(47, 4).Value2 = 2025
$ws.Cells.Item(47, 5).Value2 = -0.1533081888441812

$ws.Cells.Item(48, 1).Value2 = 45534
$ws.Cells.Item(48, 2).Value2 = 2024
$ws.Cells.Item(48, 3).Value2 = -0.05499271238530445
$ws.Cells.Item(48, 4).Value2 = 2025
$ws.Cells.Item(48, 5).Value2 = -0.1365897193907339

$ws.Cells.Item(49, 1).Value2 = 45618
$ws.Cells.Item(49, 2).Value2 = 2024
$ws.Cells.Item(49, 3).Value2 = -0.05499271238530445
$ws.Cells.Item(49, 4).Value2 = 2025
$ws.Cells.Item(49, 5).Value2 = -0.01934819856548309

$ws.Cells.Item(50, 1).Value2 = 45713
$ws.Cells.Item(50, 2).Value2 = 2025
$ws.Cells.Item(50, 3).Value2 = 0.4236015715998187
$ws.Cells.Item(50, 4).Value2 = 2026
$ws.Cells.Item(50, 5).Value2 = 0.05974827491124213

$ws.Cells.Item(51, 1).Value2 = 45800
$ws.Cells.Item(51, 2).Value2 = 2025
$ws.Cells.Item(51, 3).Value2 = 0.4305325812036687
$ws.Cells.Item(51, 4).Value2 = 2026
$ws.Cells.Item(51, 5).Value2 = 0.2136583044595852

$ws.Cells.Item(52, 1).Value2 = 45891
$ws.Cells.Item(52, 2).Value2 = 2025
$ws.Cells.Item(52, 3).Value2 = 0.4335297397760618
$ws.Cells.Item(52, 4).Value2 = 2026
$ws.Cells.Item(52, 5).Value2 = 0.2794570629465865
